$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Hays","Topeka",217,216,5474,5,8,4332,5,44,3343,5,27,8008,5,5,109805),
    @("Emporia","Viola",128,127,5114,6,36,4709,8,22,2940,8,149,8801,10,14,107859),
    @("Emporia","Junction_City",469,468,781,6,37,811,8,11,434,6,86,828,6,6,675),
    @("Greensburg","Emporia",261,260,662,6,25,501,14,25,540,6,67,534,6,6,624),
    @("Augusta","Hillsboro",22,21,641,5,17,353,5,4,182,5,33,500,5,5,385),
    @("Wichita","Andover",30,29,151,3,3,187,6,5,178,3,8,1878,6,5,389),
    @("Abilene","Leon",12,11,180,6,21,508,8,34,599,8,93,871,8,7,581)
)

$rowIndex = 3
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

$ws.Range("P10").Select()

$wb.Save()
